$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "[-, 'MCT-2A-Acionamentos Elétricos', -, -]"
$ws.Range("E2").Value = "-"

# Row 3
$ws.Range("C3").Value = "['MCT-3A-Lab. Máquinas Elétricas', -, -]"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "[-, -, 'MCT-3A-Automação Industrial', -]"

# Row 4
$ws.Range("C4").Value = "['MCT-3A-Lab. Máquinas Elétricas', -, -]"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "[-, -, 'MCT-3A-Automação Industrial', -]"

# Row 6
$ws.Range("C6").Value = "[-, -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F6").Value = "[-, -, 'MCT-3A-Automação Industrial', -]"

# Row 7
$ws.Range("C7").Value = "[-, -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F7").Value = "[-, -, 'MCT-3A-Automação Industrial', -]"

# Row 8
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "[-, -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("E8").Value = "-"

# Row 18
$ws.Range("D18").Value = "-"

# Row 19
$ws.Range("B19").Value = "-"

# Row 20
$ws.Range("B20").Value = "-"

# Row 21
$ws.Range("B21").Value = "-"
